$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 24 de Septiembre de 2020 a las 09:57"

# --- Row 7: Rusia ---
$ws.Range("B7").Value = 1128836
$ws.Range("C7").Value = 6595
$ws.Range("D7").Value = 929829
$ws.Range("E7").Value = 179059
$ws.Range("G7").Value = 149
$ws.Range("H7").Value = 19948

# --- Row 58: Singapur ---
$ws.Range("B58").Value = 57654
$ws.Range("C58").Value = 15
$ws.Range("E58").Value = 336

# --- Row 61: Suiza ---
$ws.Range("D61").Value = 42300
$ws.Range("E61").Value = 6741

# --- Rows 63/64: Moldavia/Armenia swap places (by total cases) ---
$ws.Range("A63").Value = "Armenia"
$ws.Range("B63").Value = 48251
$ws.Range("C63").Value = 374
$ws.Range("D63").Value = 43266
$ws.Range("E63").Value = 4040
$ws.Range("G63").Value = 3
$ws.Range("H63").Value = 945

$ws.Range("A64").Value = "Moldavia"
$ws.Range("B64").Value = 48232
$ws.Range("C64").Value = 0
$ws.Range("D64").Value = 36071
$ws.Range("E64").Value = 10917
$ws.Range("G64").Value = 0
$ws.Range("H64").Value = 1244

# --- Row 69: Afganistan ---
$ws.Range("B69").Value = 39170
$ws.Range("C69").Value = 25
$ws.Range("D69").Value = 32619
$ws.Range("E69").Value = 5100
$ws.Range("G69").Value = 5
$ws.Range("H69").Value = 1451

# --- Row 78: Australia ---
$ws.Range("D78").Value = 24448
$ws.Range("E78").Value = 1671

# --- Rows 130/131: Angola/Georgia swap places ---
$ws.Range("A130").Value = "Georgia"
$ws.Range("B130").Value = 4399
$ws.Range("C130").Value = 259
$ws.Range("D130").Value = 1705
$ws.Range("E130").Value = 2669
$ws.Range("H130").Value = 25

$ws.Range("A131").Value = "Angola"
$ws.Range("B131").Value = 4363
$ws.Range("C131").Value = 0
$ws.Range("D131").Value = 1473
$ws.Range("E131").Value = 2731
$ws.Range("H131").Value = 159

# --- Row 133: Lituania ---
$ws.Range("B133").Value = 4070
$ws.Range("C133").Value = 138
$ws.Range("D133").Value = 2253
$ws.Range("E133").Value = 1728
$ws.Range("G133").Value = 2
$ws.Range("H133").Value = 89

# --- Rows 143/144: Mali/Estonia swap places ---
$ws.Range("A143").Value = "Estonia"
$ws.Range("B143").Value = 3076
$ws.Range("C143").Value = 44
$ws.Range("D143").Value = 2395
$ws.Range("E143").Value = 617
$ws.Range("H143").Value = 64

$ws.Range("A144").Value = "Mali"
$ws.Range("B144").Value = 3034
$ws.Range("C144").Value = 0
$ws.Range("D144").Value = 2382
$ws.Range("E144").Value = 522
$ws.Range("H144").Value = 130

# --- Rows 214/215: Montserrat/Islas Malvinas swap places ---
$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 13
$ws.Range("E214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 0

$ws.Range("A215").Value = "Montserrat"
$ws.Range("B215").Value = 13
$ws.Range("C215").Value = 0
$ws.Range("D215").Value = 12
$ws.Range("E215").Value = 0
$ws.Range("G215").Value = 0
$ws.Range("H215").Value = 1
